$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mentioned_in_text")
$ws.Activate()

# Append the new "Supplementary Table S3.6" row describing the PRC-like
# metacell ranking process to the list of figures/tables.
$ws.Range("A28").Value = "Supplementary Table S3.6"
$ws.Range("B28").Value = "Online Supplementary Material"
$ws.Range("C28").Value = "Ranking process for choosing PRC-like metacells."

$ws.Range("C29").Select()
